$wb = $excel.ActiveWorkbook

# --- EURIBOR1M sheet: add helper formula in D2, update selection ---
$ws1 = $wb.Worksheets.Item("EURIBOR1M")
$ws1.Range("D2").Formula = "=B2*0.01"
# the formula inherited the style of its precedent cell (B2); reset to Normal
# so the cell stays unstyled, matching the source edit.
$ws1.Range("D2").Style = "Normal"
$ws1.Range("E4").Select()

# --- USDLIBOR3M sheet: rescale rates from fractional to percentage points ---
$ws4 = $wb.Worksheets.Item("USDLIBOR3M")
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws4.Cells.Item($r, 2)
    $cell.Value = $cell.Value() * 100
}

# Activate USDLIBOR3M last so it becomes the selected tab (workbook activeTab),
# which also clears tabSelected on the previously active sheet (EURIBOR6M).
$ws4.Activate()
$ws4.Range("O15").Select()
